$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.529.07"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.602.47"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.69"
$ws.Range("E5").Value = "  +3.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.58"
$ws.Range("E6").Value = "  +1.72%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.065.20"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.465.22"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.74"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.596.59"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.59"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.38"
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.09"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.34"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.29"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.993"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.23"
$ws.Range("E27").Value = "  +2.88%  "
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +6.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.85"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.69"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.845"
$ws.Range("E36").Value = "  +4.15%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "273.51"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.53"
$ws.Range("E46").Value = "  +3.58%  "
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.940.47"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.96"
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("E51").Value = "  +1.97%  "
